# Weekly fruit/vegetable price update:
# insert one new (more recent) weekly record as row 4, pushing the
# existing rows 4-36 down to 5-37.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Femacal de La Calera"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44515
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 100112022
$ws.Range("G4").Value = "Arveja Verde"
$ws.Range("H4").Value = "Perfection"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 73
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 17000
$ws.Range("M4").Value = 16521
$ws.Range("N4").Value = '$/saco 25 kilos'
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 661
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"
